# Delete row 61 (the data row whose particle index was 66/66), shifting all
# subsequent rows up by one. This matches the diff: dimension shrinks from
# A1:J151 to A1:J150, and every row from the old row 62 onward becomes the
# new row one less (61, 62, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(61).Delete() | Out-Null
